$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1 / index 1) and "全部类型" (sheet4 / index 4) both contain
# the same table of events; update the "想去人数" (column F) counts that changed.
$sheetIndexes = @(1, 4)

foreach ($idx in $sheetIndexes) {
    $ws = $wb.Worksheets.Item($idx)

    $ws.Range("F7").Value = 11958
    $ws.Range("F11").Value = 410
    $ws.Range("F13").Value = 854
    $ws.Range("F14").Value = 13468
    $ws.Range("F15").Value = 13444
    $ws.Range("F23").Value = 165
    $ws.Range("F24").Value = 171
}
